$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF (column F) values to reflect repulled data / recalculated means
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = -9
$ws.Range("F13").Value = 1
$ws.Range("F16").Value = -5
$ws.Range("F17").Value = -3
$ws.Range("F21").Value = 0
$ws.Range("F26").Value = -2
